$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Good Morning" text in E8 is replaced with "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Update the active cell / selection to E8
$ws.Range("E8").Select()
